$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 2315.074
$ws.Cells.Item(19, 9).Value = 1304.6428
$ws.Cells.Item(19, 10).Value = 3403.2307
$ws.Cells.Item(19, 11).Value = 1304.6428
$ws.Cells.Item(19, 12).Value = 3403.2307
$ws.Cells.Item(19, 13).Value = -1129.6428
$ws.Cells.Item(19, 14).Value = -3753.2307
$ws.Cells.Item(33, 8).Value = 21285.7
$ws.Cells.Item(33, 9).Value = 35120
$ws.Cells.Item(33, 10).Value = 534.25
$ws.Cells.Item(33, 11).Value = 35120
$ws.Cells.Item(33, 12).Value = 534.25
$ws.Cells.Item(33, 13).Value = -34891
$ws.Cells.Item(33, 14).Value = -992.25
$ws.Cells.Item(96, 8).Value = 3582.6
$ws.Cells.Item(96, 9).Value = 2637.6667
$ws.Cells.Item(96, 11).Value = 7913.000100000001
$ws.Cells.Item(96, 13).Value = -6540.000100000001
$ws.Cells.Item(100, 8).Value = 4949.75
$ws.Cells.Item(100, 9).Value = 5000
$ws.Cells.Item(100, 11).Value = 5000
$ws.Cells.Item(100, 13).Value = -4459
$ws.Cells.Item(101, 8).Value = 1480.375
$ws.Cells.Item(101, 10).Value = 1999.5
$ws.Cells.Item(101, 12).Value = 5998.5
$ws.Cells.Item(101, 14).Value = -9242.5
$ws.Cells.Item(111, 8).Value = 2598.6191
$ws.Cells.Item(111, 10).Value = 2329.6667
$ws.Cells.Item(111, 12).Value = 6989.000100000001
$ws.Cells.Item(111, 14).Value = -13123.0001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9186.982
$ws.Cells.Item(32, 9).Value = 6525.0835
$ws.Cells.Item(32, 11).Value = 6525.0835
$ws.Cells.Item(32, 13).Value = -6238.0835
$ws.Cells.Item(97, 8).Value = 1145.7667
$ws.Cells.Item(97, 9).Value = 1113.4445
$ws.Cells.Item(97, 11).Value = 1113.4445
$ws.Cells.Item(97, 13).Value = -617.4445000000001
$ws.Cells.Item(110, 8).Value = 29626.846
$ws.Cells.Item(110, 9).Value = 40951
$ws.Cells.Item(110, 11).Value = 40951
$ws.Cells.Item(110, 13).Value = -38906
$ws.Cells.Item(132, 8).Value = 38525816
$ws.Cells.Item(132, 9).Value = 11647.782
$ws.Cells.Item(132, 11).Value = 34943.346
$ws.Cells.Item(132, 13).Value = -32413.346

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 15462
$ws.Cells.Item(86, 9).Value = 7919.515
$ws.Cells.Item(86, 10).Value = 31018.375
$ws.Cells.Item(86, 11).Value = 7919.515
$ws.Cells.Item(86, 12).Value = 31018.375
$ws.Cells.Item(86, 13).Value = -6796.515
$ws.Cells.Item(86, 14).Value = -33264.375
$ws.Cells.Item(89, 8).Value = 15462
$ws.Cells.Item(89, 9).Value = 7919.515
$ws.Cells.Item(89, 10).Value = 31018.375
$ws.Cells.Item(89, 11).Value = 39597.575
$ws.Cells.Item(89, 12).Value = 155091.875
$ws.Cells.Item(89, 13).Value = -33981.575
$ws.Cells.Item(89, 14).Value = -166323.875
$ws.Cells.Item(134, 8).Value = 2328.3076
$ws.Cells.Item(134, 9).Value = 2015.2667
$ws.Cells.Item(134, 11).Value = 6045.800099999999
$ws.Cells.Item(134, 13).Value = -3510.800099999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 8078.5386
$ws.Cells.Item(22, 9).Value = 8078.5386
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 8078.5386
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -7728.5386
$ws.Cells.Item(31, 8).Value = 3714.718
$ws.Cells.Item(31, 9).Value = 2773.6316
$ws.Cells.Item(31, 10).Value = 4608.75
$ws.Cells.Item(31, 11).Value = 2773.6316
$ws.Cells.Item(31, 12).Value = 4608.75
$ws.Cells.Item(31, 13).Value = -2478.6316
$ws.Cells.Item(31, 14).Value = -5198.75
$ws.Cells.Item(34, 8).Value = 3714.718
$ws.Cells.Item(34, 9).Value = 2773.6316
$ws.Cells.Item(34, 10).Value = 4608.75
$ws.Cells.Item(34, 11).Value = 2773.6316
$ws.Cells.Item(34, 12).Value = 4608.75
$ws.Cells.Item(34, 13).Value = -2571.6316
$ws.Cells.Item(34, 14).Value = -5012.75
$ws.Cells.Item(58, 8).Value = 2960.9524
$ws.Cells.Item(58, 9).Value = 2829.6667
$ws.Cells.Item(58, 10).Value = 3748.6667
$ws.Cells.Item(58, 11).Value = 2829.6667
$ws.Cells.Item(58, 12).Value = 3748.6667
$ws.Cells.Item(58, 13).Value = -2626.6667
$ws.Cells.Item(58, 14).Value = -4154.6667
$ws.Cells.Item(99, 8).Value = 13309.596
$ws.Cells.Item(99, 9).Value = 17147.572
$ws.Cells.Item(99, 11).Value = 17147.572
$ws.Cells.Item(99, 13).Value = -15649.572
$ws.Cells.Item(126, 8).Value = 13309.596
$ws.Cells.Item(126, 9).Value = 17147.572
$ws.Cells.Item(126, 11).Value = 51442.716
$ws.Cells.Item(126, 13).Value = -48972.716
$ws.Cells.Item(132, 8).Value = 46006.434
$ws.Cells.Item(132, 9).Value = 60723.793
$ws.Cells.Item(132, 11).Value = 182171.379
$ws.Cells.Item(132, 13).Value = -179641.379
$ws.Cells.Item(136, 8).Value = 2960.9524
$ws.Cells.Item(136, 9).Value = 2829.6667
$ws.Cells.Item(136, 10).Value = 3748.6667
$ws.Cells.Item(136, 11).Value = 8489.000100000001
$ws.Cells.Item(136, 12).Value = 11246.0001
$ws.Cells.Item(136, 13).Value = -5939.000100000001
$ws.Cells.Item(136, 14).Value = -16346.0001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 298.8
$ws.Cells.Item(8, 9).Value = 298.8
$ws.Cells.Item(8, 11).Value = 896.4000000000001
$ws.Cells.Item(8, 13).Value = -757.4000000000001
$ws.Cells.Item(122, 8).Value = 1404
$ws.Cells.Item(122, 10).Value = 1772.5555
$ws.Cells.Item(122, 12).Value = 15952.9995
$ws.Cells.Item(122, 14).Value = -20852.9995

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2155.3823
$ws.Cells.Item(102, 9).Value = 1217.6818
$ws.Cells.Item(102, 11).Value = 1217.6818
$ws.Cells.Item(102, 13).Value = 404.3181999999999
$ws.Cells.Item(107, 8).Value = 776.5454999999999
$ws.Cells.Item(107, 9).Value = 663.9
$ws.Cells.Item(107, 11).Value = 663.9
$ws.Cells.Item(107, 13).Value = 1256.1
$ws.Cells.Item(122, 8).Value = 2678.7646
$ws.Cells.Item(122, 9).Value = 2567.2727
$ws.Cells.Item(122, 11).Value = 7701.8181
$ws.Cells.Item(122, 13).Value = -5251.8181
$ws.Cells.Item(126, 8).Value = 4834.6665
$ws.Cells.Item(126, 9).Value = 4752.25
$ws.Cells.Item(126, 11).Value = 14256.75
$ws.Cells.Item(126, 13).Value = -11786.75
$ws.Cells.Item(132, 8).Value = 2603.7
$ws.Cells.Item(132, 9).Value = 2474.9412
$ws.Cells.Item(132, 11).Value = 7424.823600000001
$ws.Cells.Item(132, 13).Value = -4894.823600000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 950
$ws.Cells.Item(22, 9).Value = 500
$ws.Cells.Item(22, 10).Value = 2300
$ws.Cells.Item(22, 11).Value = 500
$ws.Cells.Item(22, 12).Value = 2300
$ws.Cells.Item(22, 13).Value = -205
$ws.Cells.Item(22, 14).Value = -2890
$ws.Cells.Item(27, 8).Value = 950
$ws.Cells.Item(27, 9).Value = 500
$ws.Cells.Item(27, 10).Value = 2300
$ws.Cells.Item(27, 11).Value = 500
$ws.Cells.Item(27, 12).Value = 2300
$ws.Cells.Item(27, 13).Value = -393
$ws.Cells.Item(27, 14).Value = -2514
$ws.Cells.Item(46, 8).Value = 1301.9333
$ws.Cells.Item(46, 9).Value = 500.31818
$ws.Cells.Item(46, 10).Value = 3506.375
$ws.Cells.Item(46, 11).Value = 500.31818
$ws.Cells.Item(46, 12).Value = 3506.375
$ws.Cells.Item(46, 13).Value = -312.31818
$ws.Cells.Item(46, 14).Value = -3882.375
$ws.Cells.Item(82, 8).Value = 3231
$ws.Cells.Item(82, 9).Value = 2335.3333
$ws.Cells.Item(82, 11).Value = 2335.3333
$ws.Cells.Item(82, 13).Value = -1974.3333
$ws.Cells.Item(85, 8).Value = 3231
$ws.Cells.Item(85, 9).Value = 2335.3333
$ws.Cells.Item(85, 11).Value = 2335.3333
$ws.Cells.Item(85, 13).Value = -1087.3333
$ws.Cells.Item(132, 8).Value = 20340.7
$ws.Cells.Item(132, 9).Value = 22696
$ws.Cells.Item(132, 11).Value = 68088
$ws.Cells.Item(132, 13).Value = -65558

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 8000
$ws.Cells.Item(62, 10).Value = 8000
$ws.Cells.Item(62, 12).Value = 8000
$ws.Cells.Item(62, 14).Value = -9248
$ws.Cells.Item(65, 8).Value = 8000
$ws.Cells.Item(65, 10).Value = 8000
$ws.Cells.Item(65, 12).Value = 40000
$ws.Cells.Item(65, 14).Value = -46240
$ws.Cells.Item(81, 8).Value = 1059.5
$ws.Cells.Item(81, 10).Value = 1246
$ws.Cells.Item(81, 12).Value = 2492
$ws.Cells.Item(81, 14).Value = -4614
$ws.Cells.Item(84, 8).Value = 1059.5
$ws.Cells.Item(84, 10).Value = 1246
$ws.Cells.Item(84, 12).Value = 12460
$ws.Cells.Item(84, 14).Value = -23068
$ws.Cells.Item(122, 8).Value = 30770.514
$ws.Cells.Item(122, 9).Value = 35132.867
$ws.Cells.Item(122, 11).Value = 105398.601
$ws.Cells.Item(122, 13).Value = -102948.601

# ---- Special case: CRP row 22 drops the N column cell entirely ----
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCRP.Range("N22").ClearContents()

Write-Host "Applied all cell updates"